$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the naive QoQ error series matched to ifoCAST with the full-series evaluation
# Each row shifts its quarterly errors left by one column (Q0 dropped), and rows 2-14
# gain a new trailing value in column K, while rows 15-24 lose their last value (staircase shrinks).

# Row 2
$ws.Cells.Item(2, 2).Value = 9.141635797313464
$ws.Cells.Item(2, 3).Value = -8.409018285787846
$ws.Cells.Item(2, 4).Value = -0.6996216861316987
$ws.Cells.Item(2, 5).Value = 0.8832909747552129
$ws.Cells.Item(2, 6).Value = -2.082996017594043
$ws.Cells.Item(2, 7).Value = 0.18103513707042
$ws.Cells.Item(2, 8).Value = -0.2701190178563186
$ws.Cells.Item(2, 9).Value = -0.4104830813702928
$ws.Cells.Item(2, 10).Value = 0.1157559294919248
$ws.Cells.Item(2, 11).Value = -0.02607301724734923

# Row 3
$ws.Cells.Item(3, 2).Value = -8.723943454208817
$ws.Cells.Item(3, 3).Value = -1.01454685455267
$ws.Cells.Item(3, 4).Value = 0.5683658063342414
$ws.Cells.Item(3, 5).Value = -2.397921186015015
$ws.Cells.Item(3, 6).Value = -0.1338900313505515
$ws.Cells.Item(3, 7).Value = -0.5850441862772902
$ws.Cells.Item(3, 8).Value = -0.7254082497912643
$ws.Cells.Item(3, 9).Value = -0.1991692389290468
$ws.Cells.Item(3, 10).Value = -0.3409981856683208
$ws.Cells.Item(3, 11).Value = 0.3019764357337431

# Row 4
$ws.Cells.Item(4, 2).Value = -2.0709218147352
$ws.Cells.Item(4, 3).Value = -0.4880091538482882
$ws.Cells.Item(4, 4).Value = -3.454296146197544
$ws.Cells.Item(4, 5).Value = -1.190264991533081
$ws.Cells.Item(4, 6).Value = -1.64141914645982
$ws.Cells.Item(4, 7).Value = -1.781783209973794
$ws.Cells.Item(4, 8).Value = -1.255544199111576
$ws.Cells.Item(4, 9).Value = -1.39737314585085
$ws.Cells.Item(4, 10).Value = -0.7543985244487865
$ws.Cells.Item(4, 11).Value = -0.6299858181192041

# Row 5
$ws.Cells.Item(5, 2).Value = 0.3826372221226423
$ws.Cells.Item(5, 3).Value = -2.583649770226613
$ws.Cells.Item(5, 4).Value = -0.3196186155621505
$ws.Cells.Item(5, 5).Value = -0.7707727704888893
$ws.Cells.Item(5, 6).Value = -0.9111368340028634
$ws.Cells.Item(5, 7).Value = -0.3848978231406458
$ws.Cells.Item(5, 8).Value = -0.5267267698799198
$ws.Cells.Item(5, 9).Value = 0.1162478515221441
$ws.Cells.Item(5, 10).Value = 0.2406605578517265
$ws.Cells.Item(5, 11).Value = -0.1569696481858908

# Row 6
$ws.Cells.Item(6, 2).Value = -2.594300986556497
$ws.Cells.Item(6, 3).Value = -0.3302698318920341
$ws.Cells.Item(6, 4).Value = -0.7814239868187727
$ws.Cells.Item(6, 5).Value = -0.9217880503327469
$ws.Cells.Item(6, 6).Value = -0.3955490394705293
$ws.Cells.Item(6, 7).Value = -0.5373779862098034
$ws.Cells.Item(6, 8).Value = 0.1055966351922606
$ws.Cells.Item(6, 9).Value = 0.230009341521843
$ws.Cells.Item(6, 10).Value = -0.1676208645157742
$ws.Cells.Item(6, 11).Value = -0.24378957821662

# Row 7
$ws.Cells.Item(7, 2).Value = -0.4388843641081749
$ws.Cells.Item(7, 3).Value = -0.8900385190349136
$ws.Cells.Item(7, 4).Value = -1.030402582548888
$ws.Cells.Item(7, 5).Value = -0.5041635716866701
$ws.Cells.Item(7, 6).Value = -0.6459925184259441
$ws.Cells.Item(7, 7).Value = -0.003017897023880223
$ws.Cells.Item(7, 8).Value = 0.1213948093057022
$ws.Cells.Item(7, 9).Value = -0.2762353967319151
$ws.Cells.Item(7, 10).Value = -0.3524041104327608
$ws.Cells.Item(7, 11).Value = -0.07608205230983922

# Row 8
$ws.Cells.Item(8, 2).Value = -0.6819367980671746
$ws.Cells.Item(8, 3).Value = -0.8223008615811487
$ws.Cells.Item(8, 4).Value = -0.2960618507189311
$ws.Cells.Item(8, 5).Value = -0.4378907974582051
$ws.Cells.Item(8, 6).Value = 0.2050838239438588
$ws.Cells.Item(8, 7).Value = 0.3294965302734412
$ws.Cells.Item(8, 8).Value = -0.06813367576417605
$ws.Cells.Item(8, 9).Value = -0.1443023894650218
$ws.Cells.Item(8, 10).Value = 0.1320196686578998
$ws.Cells.Item(8, 11).Value = -0.0392001536538242

# Row 9
$ws.Cells.Item(9, 2).Value = -0.7058026912073482
$ws.Cells.Item(9, 3).Value = -0.1795636803451306
$ws.Cells.Item(9, 4).Value = -0.3213926270844047
$ws.Cells.Item(9, 5).Value = 0.3215819943176592
$ws.Cells.Item(9, 6).Value = 0.4459947006472416
$ws.Cells.Item(9, 7).Value = 0.04836449460962439
$ws.Cells.Item(9, 8).Value = -0.02780421909122137
$ws.Cells.Item(9, 9).Value = 0.2485178390317002
$ws.Cells.Item(9, 10).Value = 0.07729801671997623
$ws.Cells.Item(9, 11).Value = 0.2304454222287174

# Row 10
$ws.Cells.Item(10, 2).Value = -0.1340702663479956
$ws.Cells.Item(10, 3).Value = -0.2758992130872696
$ws.Cells.Item(10, 4).Value = 0.3670754083147943
$ws.Cells.Item(10, 5).Value = 0.4914881146443768
$ws.Cells.Item(10, 6).Value = 0.09385790860675949
$ws.Cells.Item(10, 7).Value = 0.01768919490591373
$ws.Cells.Item(10, 8).Value = 0.2940112530288354
$ws.Cells.Item(10, 9).Value = 0.1227914307171113
$ws.Cells.Item(10, 10).Value = 0.2759388362258526
$ws.Cells.Item(10, 11).Value = -0.001489938197266189

# Row 11
$ws.Cells.Item(11, 2).Value = -0.293453334447234
$ws.Cells.Item(11, 3).Value = 0.3495212869548299
$ws.Cells.Item(11, 4).Value = 0.4739339932844123
$ws.Cells.Item(11, 5).Value = 0.07630378724679503
$ws.Cells.Item(11, 6).Value = 0.0001350735459492769
$ws.Cells.Item(11, 7).Value = 0.2764571316688709
$ws.Cells.Item(11, 8).Value = 0.1052373093571469
$ws.Cells.Item(11, 9).Value = 0.2583847148658881
$ws.Cells.Item(11, 10).Value = -0.01904405955723064
$ws.Cells.Item(11, 11).Value = -0.182031752916177

# Row 12
$ws.Cells.Item(12, 2).Value = 0.4376286059058361
$ws.Cells.Item(12, 3).Value = 0.5620413122354185
$ws.Cells.Item(12, 4).Value = 0.1644111061978012
$ws.Cells.Item(12, 5).Value = 0.08824239249695551
$ws.Cells.Item(12, 6).Value = 0.3645644506198771
$ws.Cells.Item(12, 7).Value = 0.1933446283081531
$ws.Cells.Item(12, 8).Value = 0.3464920338168943
$ws.Cells.Item(12, 9).Value = 0.06906325939377558
$ws.Cells.Item(12, 10).Value = -0.09392443396517081
$ws.Cells.Item(12, 11).Value = -0.2180070093596886

# Row 13
$ws.Cells.Item(13, 2).Value = 0.7406159457232021
$ws.Cells.Item(13, 3).Value = 0.3429857396855849
$ws.Cells.Item(13, 4).Value = 0.2668170259847391
$ws.Cells.Item(13, 5).Value = 0.5431390841076607
$ws.Cells.Item(13, 6).Value = 0.3719192617959367
$ws.Cells.Item(13, 7).Value = 0.525066667304678
$ws.Cells.Item(13, 8).Value = 0.2476378928815592
$ws.Cells.Item(13, 9).Value = 0.0846501995226128
$ws.Cells.Item(13, 10).Value = -0.03943237587190501
$ws.Cells.Item(13, 11).Value = 0.4767206611340558

# Row 14
$ws.Cells.Item(14, 2).Value = 1.255012967438235
$ws.Cells.Item(14, 3).Value = 1.178844253737389
$ws.Cells.Item(14, 4).Value = 1.455166311860311
$ws.Cells.Item(14, 5).Value = 1.283946489548587
$ws.Cells.Item(14, 6).Value = 1.437093895057328
$ws.Cells.Item(14, 7).Value = 1.159665120634209
$ws.Cells.Item(14, 8).Value = 0.9966774272752628
$ws.Cells.Item(14, 9).Value = 0.8725948518807449
$ws.Cells.Item(14, 10).Value = 1.388747888886706
$ws.Cells.Item(14, 11).Value = 1.178844253737389

# Row 15
$ws.Cells.Item(15, 2).Value = 0.2348700177716323
$ws.Cells.Item(15, 3).Value = 0.5111920758945538
$ws.Cells.Item(15, 4).Value = 0.3399722535828299
$ws.Cells.Item(15, 5).Value = 0.4931196590915711
$ws.Cells.Item(15, 6).Value = 0.2156908846684524
$ws.Cells.Item(15, 7).Value = 0.05270319130950599
$ws.Cells.Item(15, 8).Value = -0.07137938408501182
$ws.Cells.Item(15, 9).Value = 0.444773652920949
$ws.Cells.Item(15, 10).Value = 0.2348700177716323
$ws.Cells.Item(15, 11).Value = $null

# Row 16
$ws.Cells.Item(16, 2).Value = 0.5151599734076631
$ws.Cells.Item(16, 3).Value = 0.343940151095939
$ws.Cells.Item(16, 4).Value = 0.4970875566046802
$ws.Cells.Item(16, 5).Value = 0.2196587821815615
$ws.Cells.Item(16, 6).Value = 0.0566710888226151
$ws.Cells.Item(16, 7).Value = -0.06741148657190271
$ws.Cells.Item(16, 8).Value = 0.4487415504340581
$ws.Cells.Item(16, 9).Value = 0.2388379152847414
$ws.Cells.Item(16, 10).Value = $null
$ws.Cells.Item(16, 11).Value = $null

# Row 17
$ws.Cells.Item(17, 2).Value = 0.4795802412661804
$ws.Cells.Item(17, 3).Value = 0.6327276467749217
$ws.Cells.Item(17, 4).Value = 0.3552988723518029
$ws.Cells.Item(17, 5).Value = 0.1923111789928565
$ws.Cells.Item(17, 6).Value = 0.06822860359833866
$ws.Cells.Item(17, 7).Value = 0.5843816406042994
$ws.Cells.Item(17, 8).Value = 0.3744780054549828
$ws.Cells.Item(17, 9).Value = $null
$ws.Cells.Item(17, 10).Value = $null
$ws.Cells.Item(17, 11).Value = $null

# Row 18
$ws.Cells.Item(18, 2).Value = 0.3919214649192569
$ws.Cells.Item(18, 3).Value = 0.1144926904961382
$ws.Cells.Item(18, 4).Value = -0.04849500286280822
$ws.Cells.Item(18, 5).Value = -0.172577578257326
$ws.Cells.Item(18, 6).Value = 0.3435754587486348
$ws.Cells.Item(18, 7).Value = 0.1336718235993181
$ws.Cells.Item(18, 8).Value = $null
$ws.Cells.Item(18, 9).Value = $null
$ws.Cells.Item(18, 10).Value = $null
$ws.Cells.Item(18, 11).Value = $null

# Row 19
$ws.Cells.Item(19, 2).Value = 0.0691614752440418
$ws.Cells.Item(19, 3).Value = -0.09382621811490459
$ws.Cells.Item(19, 4).Value = -0.2179087935094224
$ws.Cells.Item(19, 5).Value = 0.2982442434965384
$ws.Cells.Item(19, 6).Value = 0.08834060834722172
$ws.Cells.Item(19, 7).Value = $null
$ws.Cells.Item(19, 8).Value = $null
$ws.Cells.Item(19, 9).Value = $null
$ws.Cells.Item(19, 10).Value = $null
$ws.Cells.Item(19, 11).Value = $null

# Row 20
$ws.Cells.Item(20, 2).Value = -0.1606876400509585
$ws.Cells.Item(20, 3).Value = -0.2847702154454763
$ws.Cells.Item(20, 4).Value = 0.2313828215604846
$ws.Cells.Item(20, 5).Value = 0.02147918641116785
$ws.Cells.Item(20, 6).Value = $null
$ws.Cells.Item(20, 7).Value = $null
$ws.Cells.Item(20, 8).Value = $null
$ws.Cells.Item(20, 9).Value = $null
$ws.Cells.Item(20, 10).Value = $null
$ws.Cells.Item(20, 11).Value = $null

# Row 21
$ws.Cells.Item(21, 2).Value = -0.3143564178021929
$ws.Cells.Item(21, 3).Value = 0.201796619203768
$ws.Cells.Item(21, 4).Value = -0.00810701594554874
$ws.Cells.Item(21, 5).Value = $null
$ws.Cells.Item(21, 6).Value = $null
$ws.Cells.Item(21, 7).Value = $null
$ws.Cells.Item(21, 8).Value = $null
$ws.Cells.Item(21, 9).Value = $null
$ws.Cells.Item(21, 10).Value = $null
$ws.Cells.Item(21, 11).Value = $null

# Row 22
$ws.Cells.Item(22, 2).Value = 0.1836459624741271
$ws.Cells.Item(22, 3).Value = -0.02625767267518964
$ws.Cells.Item(22, 4).Value = $null
$ws.Cells.Item(22, 5).Value = $null
$ws.Cells.Item(22, 6).Value = $null
$ws.Cells.Item(22, 7).Value = $null
$ws.Cells.Item(22, 8).Value = $null
$ws.Cells.Item(22, 9).Value = $null
$ws.Cells.Item(22, 10).Value = $null
$ws.Cells.Item(22, 11).Value = $null

# Row 23
$ws.Cells.Item(23, 2).Value = -0.04428949692388896
$ws.Cells.Item(23, 3).Value = $null
$ws.Cells.Item(23, 4).Value = $null
$ws.Cells.Item(23, 5).Value = $null
$ws.Cells.Item(23, 6).Value = $null
$ws.Cells.Item(23, 7).Value = $null
$ws.Cells.Item(23, 8).Value = $null
$ws.Cells.Item(23, 9).Value = $null
$ws.Cells.Item(23, 10).Value = $null
$ws.Cells.Item(23, 11).Value = $null

# Row 24
$ws.Cells.Item(24, 2).Value = $null
$ws.Cells.Item(24, 3).Value = $null
$ws.Cells.Item(24, 4).Value = $null
$ws.Cells.Item(24, 5).Value = $null
$ws.Cells.Item(24, 6).Value = $null
$ws.Cells.Item(24, 7).Value = $null
$ws.Cells.Item(24, 8).Value = $null
$ws.Cells.Item(24, 9).Value = $null
$ws.Cells.Item(24, 10).Value = $null
$ws.Cells.Item(24, 11).Value = $null
